$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - "Cálculo básico" test case: input list now shows positional style [1000, 20, 1]
$ws.Range("E3").Value = "[1000, 20, 1]"

# Row 7 - "Verifica string com int/float" (capital as string) - introduce new string first
$ws.Range("D7").Value = "Verifica string com int/float"
$ws.Range("E7").Value = '["str", 20, 1]'
$ws.Range("F7").Value = 'TypeError: ("Apenas números")'

# Row 4 - "Verificar se o resultado é negativo" (capital negative)
$ws.Range("E4").Value = "[-1000, 20, 1]"
$ws.Range("F4").Value = 'ValueError: ("Apenas valores positivos")'

# Row 5 - "Verificar se o resultado é negativo" (porcentagem negative)
$ws.Range("D5").Value = "Verificar se o resultado é negativo"
$ws.Range("E5").Value = "[1000, -20, 1]"
$ws.Range("F5").Value = 'ValueError: ("Apenas valores positivos")'

# Row 6 - "Verificar se o resultado é negativo" (tempo negative)
$ws.Range("D6").Value = "Verificar se o resultado é negativo"
$ws.Range("E6").Value = "[1000, 20, -1]"
$ws.Range("F6").Value = 'ValueError: ("Apenas valores positivos")'

# Row 8 - "Verifica string com int/float" (porcentagem as string)
$ws.Range("D8").Value = "Verifica string com int/float"
$ws.Range("E8").Value = '[1000, "str", 1]'
$ws.Range("F8").Value = 'TypeError: ("Apenas números")'

# Row 9 - "Verifica string com int/float" (tempo as string)
$ws.Range("D9").Value = "Verifica string com int/float"
$ws.Range("E9").Value = '[1000, 20, "str"]'
$ws.Range("F9").Value = 'TypeError: ("Apenas números")'

# Row 10 - "Verifica se os parâmetros estão vazios"
$ws.Range("D10").Value = "Verifica se os parâmetros estão vazios"
$ws.Range("E10").Value = "[]"
$ws.Range("F10").Value = 'ValueError: ("valores vazios")'

# Update the active selection to match the saved cursor position
$ws.Range("H7").Select()
